$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) "Level proto" (A25) becomes a bold + underlined section title.
# ---------------------------------------------------------------------------
$ws.Range("A25").Font.Bold = $true
$ws.Range("A25").Font.Underline = $true

# ---------------------------------------------------------------------------
# 2) Prep the formats for the new rows up front (doesn't affect shared
#    string order, only the later value assignments do).
# ---------------------------------------------------------------------------
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A48").PasteSpecial(-4122) | Out-Null
$ws.Range("A55").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B8").Copy() | Out-Null
$ws.Range("B36").PasteSpecial(-4122) | Out-Null
$ws.Range("B50").PasteSpecial(-4122) | Out-Null
$ws.Range("B51").PasteSpecial(-4122) | Out-Null
$ws.Range("B52").PasteSpecial(-4122) | Out-Null
$ws.Range("B57").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B36").Borders.LineStyle = -4142
$ws.Range("B50").Borders.LineStyle = -4142
$ws.Range("B51").Borders.LineStyle = -4142
$ws.Range("B52").Borders.LineStyle = -4142
$ws.Range("B57").Borders.LineStyle = -4142

$ws.Range("B36").Value = 0
$ws.Range("B50").Value = 0
$ws.Range("B51").Value = 0
$ws.Range("B52").Value = 0
$ws.Range("B57").Value = 0

# ---------------------------------------------------------------------------
# 3) Text updates, in the exact order the strings were (re)typed so the
#    shared-string table is appended to in the same sequence as the source
#    edit.
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = "Ecran lobby niveau + retour + boutique"
$ws.Range("A30").Value = "Mécanisme de rangement"
$ws.Range("A29").Value = "Mécanisme de tir + physique"
$ws.Range("A36").Value = "Pause (recommencer+ retour au menu principal + boutique)"
$ws.Range("A33").Value = "Changement de scènes"
$ws.Range("A48").Value = "Level Design :"
$ws.Range("A50").Value = "Création de plusieurs niveau "
$ws.Range("A51").Value = "Courbe de difficulté"
$ws.Range("A52").Value = "Création de plusieurs monde"
$ws.Range("A41").Value = "Score de fin"
$ws.Range("A34").Value = "Scoring"
$ws.Range("A55").Value = "3D :"
$ws.Range("A57").Value = "Forme tetris"
$ws.Range("A32").Value = "Mécanisme de retry"

# ---------------------------------------------------------------------------
# 6) Cosmetic: widen column A a touch and update the saved selection /
#    scroll position to match where the author ended up editing.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 50.35

$ws.Application.Goto($ws.Range("A24"))
$ws.Range("B36").Select() | Out-Null
